$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2: C2 becomes a real number (12345678) instead of text
$ws.Range("C2").Value = 12345678

# Insert a new row 3 with the "validator_upi" record
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "validator_upi"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "12345678"
$ws.Range("D3").Value = "validator"
$ws.Range("E3").Value = "madd.hdyt@gmail.com"
$ws.Range("F3").Value = "upi"
$ws.Range("G3").Value = "dosen"
$ws.Range("H3").Value = "scopus.com"
$ws.Range("I3").Value = "scopus.com"
$ws.Range("J3").Value = "scopus.com"
$ws.Range("K3").Value = "validator"
